# Apply scheduled runner updates to Golem_Profits workbook (ALC..LTW sheets)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3149.25
$ws.Range("J17").Value = 3199
$ws.Range("L17").Value = 9597
$ws.Range("N17").Value = -9933
$ws.Range("H18").Value = 999
$ws.Range("I18").Value = 999
$ws.Range("K18").Value = 999
$ws.Range("M18").Value = -715
$ws.Range("H64").Value = 6999
$ws.Range("J64").Value = 6999
$ws.Range("L64").Value = 6999
$ws.Range("N64").Value = -7495
$ws.Range("H67").Value = 6999
$ws.Range("J67").Value = 6999
$ws.Range("L67").Value = 6999
$ws.Range("N67").Value = -8715
$ws.Range("H74").Value = 5450
$ws.Range("I74").Value = 900
$ws.Range("K74").Value = 900
$ws.Range("M74").Value = 36
$ws.Range("H77").Value = 5450
$ws.Range("I77").Value = 900
$ws.Range("K77").Value = 4500
$ws.Range("M77").Value = 180
$ws.Range("H86").Value = 2448
$ws.Range("I86").Value = 2448
$ws.Range("K86").Value = 2448
$ws.Range("M86").Value = -1325
$ws.Range("H87").Value = 73430.164
$ws.Range("J87").Value = 100151.5
$ws.Range("L87").Value = 100151.5
$ws.Range("N87").Value = -102647.5
$ws.Range("H89").Value = 2448
$ws.Range("I89").Value = 2448
$ws.Range("K89").Value = 12240
$ws.Range("M89").Value = -6624
$ws.Range("H90").Value = 73430.164
$ws.Range("J90").Value = 100151.5
$ws.Range("L90").Value = 300454.5
$ws.Range("N90").Value = -312934.5
$ws.Range("H106").Value = 3374.75
$ws.Range("I106").Value = 3374.75
$ws.Range("K106").Value = 3374.75
$ws.Range("M106").Value = -2743.75
$ws.Range("H107").Value = 2922.9167
$ws.Range("I107").Value = 3407.6
$ws.Range("K107").Value = 3407.6
$ws.Range("M107").Value = -1487.6
$ws.Range("H116").Value = 4250
$ws.Range("I116").Value = 4250
$ws.Range("K116").Value = 4250
$ws.Range("M116").Value = -808
$ws.Range("H132").Value = 2141.5
$ws.Range("I132").Value = 1669.8
$ws.Range("K132").Value = 5009.4
$ws.Range("M132").Value = -2479.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3460.111
$ws.Range("I32").Value = 3460.111
$ws.Range("K32").Value = 3460.111
$ws.Range("M32").Value = -3173.111
$ws.Range("H61").Value = 6923.5
$ws.Range("I61").Value = 6923.5
$ws.Range("K61").Value = 6923.5
$ws.Range("M61").Value = -6711.5
$ws.Range("H132").Value = 2028.4546
$ws.Range("I132").Value = 1231.3
$ws.Range("K132").Value = 3693.9
$ws.Range("M132").Value = -1163.9
$ws.Range("H136").Value = 6923.5
$ws.Range("I136").Value = 6923.5
$ws.Range("K136").Value = 20770.5
$ws.Range("M136").Value = -18220.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1926
$ws.Range("I86").Value = 1687.4445
$ws.Range("J86").Value = 2999.5
$ws.Range("K86").Value = 1687.4445
$ws.Range("L86").Value = 2999.5
$ws.Range("M86").Value = -564.4445000000001
$ws.Range("N86").Value = -5245.5
$ws.Range("H89").Value = 1926
$ws.Range("I89").Value = 1687.4445
$ws.Range("J89").Value = 2999.5
$ws.Range("K89").Value = 8437.2225
$ws.Range("L89").Value = 14997.5
$ws.Range("M89").Value = -2821.2225
$ws.Range("N89").Value = -26229.5
$ws.Range("H96").Value = 2733.3333
$ws.Range("I96").Value = 2733.3333
$ws.Range("K96").Value = 2733.3333
$ws.Range("M96").Value = 12.66670000000022
$ws.Range("H107").Value = 3553.5
$ws.Range("I107").Value = 3567.2
$ws.Range("K107").Value = 3567.2
$ws.Range("M107").Value = -1647.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 25833.334
$ws.Range("J93").Value = 56000
$ws.Range("L93").Value = 56000
$ws.Range("N93").Value = -59744
$ws.Range("H105").Value = 449.625
$ws.Range("J105").Value = 499
$ws.Range("L105").Value = 499
$ws.Range("N105").Value = -3993
$ws.Range("H134").Value = 1746.75
$ws.Range("I134").Value = 893.5
$ws.Range("K134").Value = 2680.5
$ws.Range("M134").Value = -145.5
$ws.Range("H141").Value = 497708.3
$ws.Range("J141").Value = 497708.3
$ws.Range("L141").Value = 497708.3
$ws.Range("N141").Value = -508068.3

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 160.57692
$ws.Range("J2").Value = 257.5
$ws.Range("L2").Value = 1545
$ws.Range("N2").Value = -1771
$ws.Range("H34").Value = 3010.4043
$ws.Range("J34").Value = 3124.1777
$ws.Range("L34").Value = 9372.5331
$ws.Range("N34").Value = -9540.5331
$ws.Range("H38").Value = 44.5
$ws.Range("I38").Value = 44.5
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 133.5
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 213.5
$ws.Range("N38").ClearContents()
$ws.Range("H39").Value = 1999
$ws.Range("J39").Value = 1999
$ws.Range("L39").Value = 5997
$ws.Range("N39").Value = -6585
$ws.Range("H55").Value = 3371.5813
$ws.Range("I55").Value = 1309.5714
$ws.Range("J55").Value = 3772.5278
$ws.Range("K55").Value = 3928.7142
$ws.Range("L55").Value = 11317.5834
$ws.Range("M55").Value = -3751.7142
$ws.Range("N55").Value = -11671.5834
$ws.Range("H86").Value = 205.8
$ws.Range("I86").Value = 198.28572
$ws.Range("J86").Value = 223.33333
$ws.Range("K86").Value = 594.85716
$ws.Range("L86").Value = 669.99999
$ws.Range("M86").Value = 591.14284
$ws.Range("N86").Value = -3041.99999
$ws.Range("H89").Value = 205.8
$ws.Range("I89").Value = 198.28572
$ws.Range("J89").Value = 223.33333
$ws.Range("K89").Value = 1784.57148
$ws.Range("L89").Value = 2009.99997
$ws.Range("M89").Value = 4143.428519999999
$ws.Range("N89").Value = -13865.99997
$ws.Range("H131").Value = 1798.3334
$ws.Range("J131").Value = 1798.3334
$ws.Range("L131").Value = 5395.0002
$ws.Range("N131").Value = -15475.0002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 125.695656
$ws.Range("J2").Value = 116.42857
$ws.Range("L2").Value = 116.42857
$ws.Range("N2").Value = -342.42857
$ws.Range("H46").Value = 2499.75
$ws.Range("I46").Value = 2499.75
$ws.Range("K46").Value = 2499.75
$ws.Range("M46").Value = -2343.75
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 1631.5
$ws.Range("I80").Value = 1268.8
$ws.Range("J80").Value = 2236
$ws.Range("K80").Value = 1268.8
$ws.Range("L80").Value = 2236
$ws.Range("M80").Value = -270.8
$ws.Range("N80").Value = -4232
$ws.Range("H83").Value = 1631.5
$ws.Range("I83").Value = 1268.8
$ws.Range("J83").Value = 2236
$ws.Range("K83").Value = 6344
$ws.Range("L83").Value = 11180
$ws.Range("M83").Value = -1352
$ws.Range("N83").Value = -21164

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 960.1
$ws.Range("I22").Value = 516.5
$ws.Range("J22").Value = 1625.5
$ws.Range("K22").Value = 516.5
$ws.Range("L22").Value = 1625.5
$ws.Range("M22").Value = -221.5
$ws.Range("N22").Value = -2215.5
$ws.Range("H27").Value = 960.1
$ws.Range("I27").Value = 516.5
$ws.Range("J27").Value = 1625.5
$ws.Range("K27").Value = 516.5
$ws.Range("L27").Value = 1625.5
$ws.Range("M27").Value = -409.5
$ws.Range("N27").Value = -1839.5
$ws.Range("H46").Value = 916.2105
$ws.Range("J46").Value = 991
$ws.Range("L46").Value = 991
$ws.Range("N46").Value = -1367
$ws.Range("H132").Value = 3723.5
$ws.Range("I132").Value = 3479.8
$ws.Range("J132").Value = 3897.5715
$ws.Range("K132").Value = 10439.4
$ws.Range("L132").Value = 11692.7145
$ws.Range("M132").Value = -7909.400000000001
$ws.Range("N132").Value = -16752.7145

